# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing header style used by the other header cells, and populate
# the data rows 2-42 with the corresponding values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - copy the style of the existing header cell H1 so the
# new headers (I1, J1) match formatting (bold, centered, bordered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()

# Data values for I2:J42
$data = @(
    @(5, 6),
    @(10, 10),
    @(1, 3),
    @(6, 7),
    @(1, 3),
    @(2, 4),
    @(5, 6),
    @(7, 8),
    @(6, 7),
    @(8, 9),
    @(6, 7),
    @(8, 9),
    @(5, 6),
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(1, 5),
    @(8, 8),
    @(1, 2),
    @(9, 9),
    @(11, 11),
    @(1, 3),
    @(10, 10),
    @(5, 7),
    @(2, 4),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(7, 7),
    @(7, 8),
    @(6, 8),
    @(5, 7),
    @(7, 8),
    @(6, 7),
    @(7, 8),
    @(3, 5),
    @(5, 8),
    @(3, 6),
    @(6, 8),
    @(6, 8),
    @(1, 2)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $ws.Cells.Item($rowNum, 9).Value = $data[$r][0]
    $ws.Cells.Item($rowNum, 10).Value = $data[$r][1]
}
